# UC011 - Listar Solicitações de Diárias
# Version bump 1.0 -> 1.2.5, revision metadata update, and a handful of
# wording / numbering fixes in the use-case body.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $raw = $p.Range.Text
    $text = $raw.TrimEnd("`r", "`a")

    # --- Revision history table (single-run paragraphs) ---
    if ($text.Equals("1.0")) {
        $p.Range.Text = "1.2.5"
    }
    elseif ($text.Equals("Creation")) {
        $p.Range.Text = "Update"
    }
    elseif ($text.Equals("Fabrício Araújo")) {
        $p.Range.Text = "Julio Paiva"
    }
    elseif ($text.Equals("09/11/2020")) {
        $p.Range.Text = "31/05/2023"
    }

    # --- Preconditions table cell: fix typo + add final period ---
    elseif ($text.Equals("O usuario devidamente autenticado e na tela inicial do sistema")) {
        $p.Range.Text = "O usuário devidamente autenticado e na tela inicial do sistema."
    }

    # --- Basic flow step renumbering (multi-run / bookmarked paragraphs) ---
    elseif ($text.Contains("4. Chefe Visualiza os registros de solicitações de diária. af[1,2,3,4]")) {
        $p.Range.Find.Execute("4. Chefe Visualiza os registros de solicitações de diária. af[1,2,3,4]", `
                               $true, $false, $false, $false, $false, `
                               $true, 1, $false, `
                               "3. Chefe Visualiza os registros de solicitações de diária. af[1,2,3,4]", 2) | Out-Null
    }
    elseif ($text.Contains("5. System Exibe os registros de solicitações de diária. ")) {
        $p.Range.Find.Execute("5. System Exibe os registros de solicitações de diária. ", `
                               $true, $false, $false, $false, $false, `
                               $true, 1, $false, `
                               "4. System Exibe os registros de solicitações de diária. ", 2) | Out-Null
    }

    # --- AF[1] step 2: add missing period at end of sentence ---
    elseif ($text.Contains("Apresenta a tela de Detalhar Diárias ")) {
        $p.Range.Find.Execute("Apresenta a tela de Detalhar Diárias ", `
                               $true, $false, $false, $false, $false, `
                               $true, 1, $false, `
                               "Apresenta a tela de Detalhar Diárias. ", 2) | Out-Null
    }

    # --- AF[3] step renumbering: 4 -> 2 ---
    elseif ($text.Contains("4. System Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário. ")) {
        $p.Range.Find.Execute("4. System Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário. ", `
                               $true, $false, $false, $false, $false, `
                               $true, 1, $false, `
                               "2. System Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário. ", 2) | Out-Null
    }
}
